$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value2 = '26.325.13'
$ws.Cells.Item(2, 5).Value2 = '  +1.88%  '
$ws.Cells.Item(3, 4).Value2 = '1.649.08'
$ws.Cells.Item(3, 5).Value2 = '  +0.46%  '
$ws.Cells.Item(4, 5).Value2 = '  -0.29%  '
$ws.Cells.Item(5, 4).Value2 = '''217.58'
$ws.Cells.Item(5, 5).Value2 = '  +0.70%  '
$ws.Cells.Item(6, 4).Value2 = '''0.508'
$ws.Cells.Item(6, 5).Value2 = '  +0.43%  '
$ws.Cells.Item(7, 5).Value2 = '  -0.34%  '
$ws.Cells.Item(8, 4).Value2 = '''0.259'
$ws.Cells.Item(8, 5).Value2 = '  +0.18%  '
$ws.Cells.Item(10, 4).Value2 = '''20.07'
$ws.Cells.Item(10, 5).Value2 = '  +1.65%  '
$ws.Cells.Item(11, 5).Value2 = '  -0.12%  '
$ws.Cells.Item(12, 2).Value2 = 'Polkadot'
$ws.Cells.Item(12, 3).Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(12, 4).Value2 = '''4.31'
$ws.Cells.Item(12, 5).Value2 = '  +0.87%  '
$ws.Cells.Item(13, 2).Value2 = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(13, 3).Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(13, 4).Value2 = '1.874.77'
$ws.Cells.Item(13, 5).Value2 = '  +0.38%  '
$ws.Cells.Item(14, 2).Value2 = 'WrappedEther'
$ws.Cells.Item(14, 3).Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value2 = '1.669.13'
$ws.Cells.Item(14, 5).Value2 = '  +1.70%  '
$ws.Cells.Item(15, 4).Value2 = '''0.552'
$ws.Cells.Item(15, 5).Value2 = '  -1.90%  '
$ws.Cells.Item(16, 4).Value2 = '0.0₃0768'
$ws.Cells.Item(16, 5).Value2 = '  -0.08%  '
$ws.Cells.Item(17, 4).Value2 = '''63.61'
$ws.Cells.Item(17, 5).Value2 = '  +0.79%  '
$ws.Cells.Item(18, 4).Value2 = '26.299.65'
$ws.Cells.Item(18, 5).Value2 = '  +1.61%  '
$ws.Cells.Item(19, 5).Value2 = '  -0.21%  '
$ws.Cells.Item(20, 4).Value2 = '''196.97'
$ws.Cells.Item(20, 5).Value2 = '  +2.03%  '
$ws.Cells.Item(21, 4).Value2 = '''4.46'
$ws.Cells.Item(21, 5).Value2 = '  -0.37%  '
$ws.Cells.Item(22, 4).Value2 = '''10.08'
$ws.Cells.Item(22, 5).Value2 = '  +0.90%  '
$ws.Cells.Item(23, 4).Value2 = '''6.36'
$ws.Cells.Item(23, 5).Value2 = '  +0.03%  '
$ws.Cells.Item(24, 5).Value2 = '  -2.61%  '
$ws.Cells.Item(25, 4).Value2 = '''143.20'
$ws.Cells.Item(25, 5).Value2 = '  +0.80%  '
$ws.Cells.Item(26, 4).Value2 = '''0.999'
$ws.Cells.Item(26, 5).Value2 = '  -0.33%  '
$ws.Cells.Item(27, 5).Value2 = '  +2.09%  '
$ws.Cells.Item(28, 5).Value2 = '  +0.36%  '
$ws.Cells.Item(29, 4).Value2 = '''15.67'
$ws.Cells.Item(29, 5).Value2 = '  +0.72%  '
$ws.Cells.Item(30, 4).Value2 = '''1.25'
$ws.Cells.Item(30, 5).Value2 = '  +1.27%  '
$ws.Cells.Item(31, 4).Value2 = '''0.0510'
$ws.Cells.Item(31, 5).Value2 = '  +2.71%  '
$ws.Cells.Item(32, 4).Value2 = '''3.36'
$ws.Cells.Item(32, 5).Value2 = '  +0.38%  '
$ws.Cells.Item(33, 5).Value2 = '  +0.12%  '
$ws.Cells.Item(34, 5).Value2 = '  +1.88%  '
$ws.Cells.Item(35, 5).Value2 = '  +0.88%  '
$ws.Cells.Item(36, 4).Value2 = '''0.917'
$ws.Cells.Item(36, 5).Value2 = '  +0.77%  '
$ws.Cells.Item(37, 4).Value2 = '1.142.42'
$ws.Cells.Item(38, 5).Value2 = '  +1.72%  '
$ws.Cells.Item(39, 4).Value2 = '''2.49'
$ws.Cells.Item(39, 5).Value2 = '  -1.38%  '
$ws.Cells.Item(41, 4).Value2 = '''0.998'
$ws.Cells.Item(41, 5).Value2 = '  -0.26%  '
$ws.Cells.Item(42, 4).Value2 = '''5.68'
$ws.Cells.Item(42, 5).Value2 = '  +2.17%  '
$ws.Cells.Item(43, 4).Value2 = '''100.42'
$ws.Cells.Item(43, 5).Value2 = '  -0.41%  '
$ws.Cells.Item(44, 4).Value2 = '''0.803'
$ws.Cells.Item(44, 5).Value2 = '  -0.62%  '
$ws.Cells.Item(45, 4).Value2 = '1.783.72'
$ws.Cells.Item(45, 5).Value2 = '  +0.39%  '
$ws.Cells.Item(46, 2).Value2 = 'BabyDogeCoin'
$ws.Cells.Item(46, 3).Value2 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(46, 4).Value2 = '0.0₆0110'
$ws.Cells.Item(46, 5).Value2 = '  -2.20%  '
$ws.Cells.Item(47, 2).Value2 = 'Aave'
$ws.Cells.Item(47, 3).Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).Value2 = '''56.49'
$ws.Cells.Item(47, 5).Value2 = '  +1.83%  '
$ws.Cells.Item(48, 2).Value2 = 'RenderToken'
$ws.Cells.Item(48, 3).Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(48, 4).Value2 = '''1.48'
$ws.Cells.Item(48, 5).Value2 = '  +2.60%  '
$ws.Cells.Item(49, 2).Value2 = 'Cronos'
$ws.Cells.Item(49, 3).Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).Value2 = '''0.0516'
$ws.Cells.Item(49, 5).Value2 = '  +2.27%  '
$ws.Cells.Item(50, 2).Value2 = 'EnergySwap'
$ws.Cells.Item(50, 3).Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value2 = '''7.74'
$ws.Cells.Item(50, 5).Value2 = '  +3.21%  '
$ws.Cells.Item(51, 2).Value2 = 'Mantle'
$ws.Cells.Item(51, 3).Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51, 4).Value2 = '''0.417'
$ws.Cells.Item(51, 5).Value2 = '  -0.38%  '

# Reset quote-prefix style introduced by forcing numeric-looking text values,
# so the cells keep the original 'General'/default style (no visible change).
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(51, 4).Style = "Normal"
